$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80
$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value = "Maule"
$ws.Cells.Item($row, 4).Value = 44595
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(79, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Black Amber"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 300
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 9000
$ws.Cells.Item($row, 16).Value = 9000
$ws.Cells.Item($row, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 500
$ws.Cells.Item($row, 20).Value = 18
